$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 data values (daily COVID figures, including Wales improvements)
$ws.Range("A2").Value = 43923
$ws.Range("B2").Value = 33718
$ws.Range("C2").Value = 4244
$ws.Range("D2").Value = 2921
$ws.Range("E2").Value = 569
$ws.Range("F2").Value = 28221
$ws.Range("G2").Value = 2698
$ws.Range("H2").Value = 2602
$ws.Range("I2").Value = 76
$ws.Range("J2").Value = 2121
$ws.Range("K2").Value = 117
$ws.Range("L2").Value = 774
$ws.Range("M2").Value = 30

# Remove the (now unused) shaded fill from the date cell's style, and drop
# the duplicate "165" number-format style that only differed by fill.
$ws.Range("A2").Interior.Pattern = -4142  # xlPatternNone
$ws.Range("G2").Style = $ws.Range("F2").Style

# Cosmetic formatting refresh (row heights / default gridline spacing)
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Columns("A:H").ColumnWidth = 14.43
